$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) text - order of assignment matters for shared string table order
$ws.Range("E1").Value = "Three"
$ws.Range("C1").Value = "Add"
$ws.Range("D1").Value = "test"

# Add new row 4 with data, mirroring row 3
$ws.Range("A4").Value = 45
$ws.Range("B4").Value = 565
$ws.Range("C4").Value = 6767
$ws.Range("D4").Value = 67
$ws.Range("E4").Value = 343

# Update the active selection to D1
$ws.Range("D1").Select()
